$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for Price (column D) cells whose new values
# look numeric, so Excel keeps them as text instead of auto-converting to a number.
$priceTextRows = @(5,7,8,9,10,11,13,14,15,16,17,18,22,23,24,25,26,27,28,29,30,31,32,33,34,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $priceTextRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.304.84'
$ws.Range("E2").Value = '  +0.82%  '

$ws.Range("D3").Value = '1.823.12'
$ws.Range("E3").Value = '  -0.19%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '314.44'
$ws.Range("E5").Value = '  +0.52%  '

$ws.Range("D7").Value = '0.4492'
$ws.Range("E7").Value = '  -1.83%  '

$ws.Range("D8").Value = '0.3781'
$ws.Range("E8").Value = '  +1.44%  '

$ws.Range("D9").Value = '0.07436'
$ws.Range("E9").Value = '  +1.75%  '

$ws.Range("D10").Value = '0.8860'
$ws.Range("E10").Value = '  +2.81%  '

$ws.Range("D11").Value = '20.96'
$ws.Range("E11").Value = '  +0.26%  '

$ws.Range("D12").Value = '1.824.48'
$ws.Range("E12").Value = '  -0.11%  '

$ws.Range("D13").Value = '6.730'
$ws.Range("E13").Value = '  +0.54%  '

$ws.Range("D14").Value = '5.446'
$ws.Range("E14").Value = '  +1.76%  '

$ws.Range("D15").Value = '93.53'
$ws.Range("E15").Value = '  +0.48%  '

$ws.Range("D16").Value = '0.07122'
$ws.Range("E16").Value = '  +0.27%  '

$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("D18").Value = '0.000008803'
$ws.Range("E18").Value = '  -0.52%  '

$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("E20").Value = '  +0.87%  '

$ws.Range("D21").Value = '27.319.22'
$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("D22").Value = '5.388'
$ws.Range("E22").Value = '  +3.81%  '

$ws.Range("D23").Value = '10.95'
$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("D24").Value = '1.967'
$ws.Range("E24").Value = '  -1.73%  '

$ws.Range("D25").Value = '151.61'
$ws.Range("E25").Value = '  -0.09%  '

$ws.Range("D26").Value = '2.307'
$ws.Range("E26").Value = '  +3.24%  '

$ws.Range("D27").Value = '18.63'
$ws.Range("E27").Value = '  +0.58%  '

$ws.Range("D28").Value = '5.382'
$ws.Range("E28").Value = '  +1.95%  '

$ws.Range("D29").Value = '117.78'
$ws.Range("E29").Value = '  +0.45%  '

$ws.Range("D30").Value = '0.08907'
$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("D31").Value = '0.7903'
$ws.Range("E31").Value = '  +4.25%  '

$ws.Range("D32").Value = '1.201'
$ws.Range("E32").Value = '  +0.56%  '

$ws.Range("D33").Value = '4.603'
$ws.Range("E33").Value = '  +2.85%  '

$ws.Range("D34").Value = '2.913'
$ws.Range("E34").Value = '  -1.51%  '

$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("D36").Value = '1.110'
$ws.Range("E36").Value = '  +0.80%  '

$ws.Range("D37").Value = '0.01980'
$ws.Range("E37").Value = '  +0.39%  '

$ws.Range("D38").Value = '0.05299'
$ws.Range("E38").Value = '  +0.59%  '

$ws.Range("D39").Value = '7.383'
$ws.Range("E39").Value = '  +2.65%  '

$ws.Range("D40").Value = '0.5331'
$ws.Range("E40").Value = '  -0.53%  '

$ws.Range("D41").Value = '2.867'
$ws.Range("E41").Value = '  -0.56%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1712'
$ws.Range("E42").Value = '  -0.33%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '2.316'
$ws.Range("E43").Value = '  +17.50%  '

$ws.Range("D44").Value = '8.651'
$ws.Range("E44").Value = '  +0.76%  '

$ws.Range("D45").Value = '0.5062'
$ws.Range("E45").Value = '  -3.86%  '

$ws.Range("D46").Value = '10.65'
$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("D47").Value = '1.698'
$ws.Range("E47").Value = '  +1.03%  '

$ws.Range("D48").Value = '105.26'
$ws.Range("E48").Value = '  -0.52%  '

$ws.Range("D49").Value = '1.000'
$ws.Range("E49").Value = '  -0.07%  '

$ws.Range("D50").Value = '0.06405'
$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("D51").Value = '66.05'
$ws.Range("E51").Value = '  +3.99%  '
